# A new weekly price record was inserted as the first data row (row 13).
# This shifts every existing data row (13-79) down by one (to 14-80); the
# former last row (79) becomes the new last row (80). Inserting a whole row
# (rather than writing every cell of every row) lets Excel carry the
# existing data/formatting down automatically, matching the target file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value = 10
$ws.Cells.Item(13, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(13, 3).Value = "La Araucanía"
$ws.Cells.Item(13, 4).Value = 45069
$ws.Cells.Item(13, 5).Value = 9
$ws.Cells.Item(13, 6).Value = 100112042
$ws.Cells.Item(13, 7).Value = "Locoto"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 80
$ws.Cells.Item(13, 11).Value = 4400
$ws.Cells.Item(13, 12).Value = 4400
$ws.Cells.Item(13, 13).Value = 4400
$ws.Cells.Item(13, 14).Value = "$/kilo"
$ws.Cells.Item(13, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(13, 16).Value = 4400
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = "Hortaliza"
